$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.945
$ws.Range("A14").Value = -21.683
$ws.Range("A21").Value = -19.937
$ws.Range("D22").Value = -7.976999999999999
$ws.Range("A23").Value = -20.203
$ws.Range("D24").Value = -7.327000000000001
$ws.Range("A25").Value = -21.632
$ws.Range("A26").Value = -21.662
$ws.Range("D28").Value = -8.166
$ws.Range("A29").Value = -21.344
$ws.Range("D36").Value = -7.540999999999999
$ws.Range("D45").Value = -7.619
$ws.Range("D48").Value = -7.540999999999999
$ws.Range("D49").Value = -8.320000000000002
$ws.Range("D52").Value = -7.905000000000001
$ws.Range("A53").Value = -22.068
$ws.Range("D53").Value = -8.068999999999999
$ws.Range("D54").Value = -8.183
$ws.Range("A57").Value = -22.219
$ws.Range("A59").Value = -22.461
$ws.Range("A69").Value = -21.606
$ws.Range("D70").Value = -6.853
$ws.Range("A79").Value = -20.778
$ws.Range("A83").Value = -22.024
$ws.Range("D86").Value = -8.252000000000001
$ws.Range("D87").Value = -8.228
$ws.Range("D89").Value = -8.177000000000001
$ws.Range("A91").Value = -20.666
$ws.Range("A93").Value = -21.508
$ws.Range("D101").Value = -8.183000000000002
$ws.Range("A103").Value = -22.086
